$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows at row 23 (shifts old rows 23-27 down to 29-33)
$ws.Rows.Item(23).Resize(6).Insert()

# Update row content for rows 21-28 (column A = event id, column B = event text)
$ws.Range("A21").Value = 'e020'
$ws.Range("B21").Value = '<Bold>e020 Enemy Strength Check - Choose Area</Bold> 
<InlineUIContainer><Button Content=''r4.53'' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
Check any one adjacent area to your task force for estimating enemy strength. Click on one of the adjacent regions highlighted blue.
<LineBreak/><LineBreak/>'
$ws.Rows.Item(21).RowHeight = 75

$ws.Range("A22").Value = 'e021'
$ws.Range("B22").Value = '<Bold>e020 Enemy Strength Check Roll</Bold> 
<InlineUIContainer><Button Content=''r4.53'' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
Roll 1D and consult the <InlineUIContainer><Button Content=''Resistance'' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer> 
Table. The area is marked with a Light, Medium, or Heavy marker.
<LineBreak/><LineBreak/>
Die Roll =  <InlineUIContainer><Image Name=''DieRoll'' Height=''21'' Width=''21'' > </Image></InlineUIContainer> 
<LineBreak/><LineBreak/>'
$ws.Rows.Item(22).RowHeight = 120

$ws.Range("A23").Value = 'e022'
$ws.Range("B23").Value = '<Bold>e021 Choose Operations</Bold> 
<InlineUIContainer><Button Content=''r4.54'' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
Choose one of following options. To see the options, select the ''e###'' button. To choose the option, select the other buttons. Each option uses up time per the 
<InlineUIContainer><Button Content=''Time'' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer> 
Table.
<LineBreak/><LineBreak/>
<InlineUIContainer><Button Content=''e023'' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer>
 Perform an additional <InlineUIContainer><Button Content=''Area Check'' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer><LineBreak/>
<InlineUIContainer><Button Content=''e024'' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer> Call for Artillery 
<InlineUIContainer><Button Content=''Support'' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer><LineBreak/>
<InlineUIContainer><Button Content=''e025'' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer> Call for Air 
<InlineUIContainer><Button Content=''Strike'' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer><LineBreak/>
<InlineUIContainer><Button Content=''e026'' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer> Attempt to 
<InlineUIContainer><Button Content=''Resupply'' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer><LineBreak/>
<InlineUIContainer><Button Content=''e027'' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer> Click on highlighted area to Enter
<LineBreak/><LineBreak/>'
$ws.Rows.Item(23).RowHeight = 255

$ws.Range("A24").Value = 'e023'
$ws.Range("B24").Value = '<Bold>e023 Area Check</Bold> 
<InlineUIContainer><Button Content=''r4.54.1'' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
Check any one additional adjacent area for estimated enemy strength by rolling on the 
<InlineUIContainer><Button Content=''Resistance'' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer> Table. Click on one of the adjacent regions highlighted blue.'
$ws.Rows.Item(24).RowHeight = 75

$ws.Range("A25").Value = 'e024'
$ws.Range("B25").Value = '<Bold>e024 Call for Artillery Support</Bold> 
<InlineUIContainer><Button Content=''r4.54.2'' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
Call to hit an area adjacent to your task force. Consult the 
<InlineUIContainer><Button Content=''Time'' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer> 
to see if Artillery Support arrives. If successful, an Artillery Support marker on the area. Click on one of the adjacent regions highlighted blue.'
$ws.Rows.Item(25).RowHeight = 90

$ws.Range("A26").Value = 'e025'
$ws.Range("B26").Value = '<Bold>e025 Call for Air Strike</Bold> 
<InlineUIContainer><Button Content=''r4.54.3'' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
Call to hit an area adjacent to your task force. Consult the 
<InlineUIContainer><Button Content=''Time'' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer> 
to see if Air Strike arrives. If successful, an Air Strike marker is placed on the area.'
$ws.Rows.Item(26).RowHeight = 90

$ws.Range("A27").Value = 'e026'
$ws.Range("B27").Value = '<Bold>e026 Attempt to Resupply</Bold> 
<InlineUIContainer><Button Content=''r4.54.4'' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
Consult the 
<InlineUIContainer><Button Content=''Time'' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer> 
to see if resupply occurs. If successful, you may relead your tank with ammo. '
$ws.Rows.Item(27).RowHeight = 90

$ws.Range("A28").Value = 'e027'
$ws.Range("B28").Value = '<Bold>e027 Enter Adjacent Area</Bold> 
<InlineUIContainer><Button Content=''r4.54.5'' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>'
$ws.Rows.Item(28).RowHeight = 45

# Restore selection/view near the edited area
$ws.Range("B28").Select()